$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Comments" column (E) so that
# "SecondVaccinationDate" becomes the new column E header and "Comments"
# shifts one column to the right, into column F.
$ws.Range("E1").EntireColumn.Insert()
$ws.Range("E1").Value = "SecondVaccinationDate"
$ws.Columns("E").AutoFit()

# A couple of cells in row 2 carry left-over number formats (date format on
# C2, a minute:second format on I2/J2) even though they hold no values.
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("I2").NumberFormat = "mm:ss.0"
$ws.Range("J2").NumberFormat = "mm:ss.0"

# Move the active selection to F2, matching the saved view state.
$ws.Range("F2").Select() | Out-Null
